$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original "SCC_PROMPT" row (row 45) is being split into two checks: one
# for headphones (relabelled SCC_PROMPT_HP, body text unchanged) and a brand
# new one for loudspeakers (SCC_PROMPT_LS), appended as a new row 53.

# 1) Seed the new row's key first so the shared-string table grows in the
#    canonical order (SCC_PROMPT_LS, SCC_PROMPT_HP, German text, English text).
$ws.Cells.Item(53, 1).Value = "SCC_PROMPT_LS"

# 2) Relabel row 45's key from SCC_PROMPT -> SCC_PROMPT_HP (B45/C45 untouched).
$ws.Cells.Item(45, 1).Value = "SCC_PROMPT_HP"

# 3) Fill in the loudspeaker prompt body text (German, then English).
$ws.Cells.Item(53, 2).Value = "<p>Sie haben angegeben, dass Sie ein anderes Gerät als Lautsprecher zur Tonwiedergabe nutzen. Für diese Befragung ist die Nutzung von Lautsprechern jedoch Grundvoraussetzung.</p><p>Im weiteren Verlauf des Fragebogens wurden spezielle Höraufgaben integriert. Mit diesen Aufgaben wird überprüft, ob von den Teilnehmern wirklich Lautsprecher getragen werden.</p><p>Wir bitten Sie daher <strong>AB jetzt unbedingt Lautsprecher zu benutzen</strong>, damit Sie die Befragung erfolgreich abschließen können. Falls Sie ab jetzt keine Lautsprecher benutzen, werden Sie im Verlauf der Befragung automatisch aussortiert. In Anbetracht der Länge der Befragung wäre dies sehr ärgerlich.</p>"
$ws.Cells.Item(53, 3).Value = "<p>You indicated that you are using a device different from loudspeakers to reproduce sound. However, the use of loudspeakers is a basic requirement for this survey.</p><p>We integrated special listening tasks in the further course of the questionnaire. These tasks check whether the participants are actually using loudspeakers. </p><p>We therefore ask you <strong> to use loudspeakers from now on</strong> to successfully complete the survey. If you do not use loudspeakers from now on, you will be automatically screened out during the survey. Given the length of the questionnaire, this would be a pity.</p>"

# Match the existing table's alignment pattern (top for A/B, centered for C).
$ws.Cells.Item(53, 1).VerticalAlignment = -4160
$ws.Cells.Item(53, 2).VerticalAlignment = -4160
$ws.Cells.Item(53, 3).VerticalAlignment = -4108

# Update the active selection to reflect the new end of the table.
$ws.Range("C57").Select() | Out-Null
